# Updates the cryptocurrency price/volume snapshot in columns D (Price) and
# E (Volume(1h)) to the latest scraped values. Cell contents are plain text
# (prices use locale-style "." group separators and are not valid numeric
# literals, e.g. "60.999.22"; percentages keep their original padding, e.g.
# "  -2.06%  "), so every value is written with a leading quote-prefix to
# force text storage and the style is restored to "Normal" right after so
# the quote-prefix flag doesn't linger in the cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '60.999.22' },
    @{ Cell = 'E2'; Value = '  -2.06%  ' },
    @{ Cell = 'D3'; Value = '2.432.27' },
    @{ Cell = 'E3'; Value = '  -0.53%  ' },
    @{ Cell = 'D4'; Value = '0.997' },
    @{ Cell = 'E4'; Value = '  -0.13%  ' },
    @{ Cell = 'D5'; Value = '571.37' },
    @{ Cell = 'E5'; Value = '  -2.33%  ' },
    @{ Cell = 'D6'; Value = '140.34' },
    @{ Cell = 'E6'; Value = '  -2.52%  ' },
    @{ Cell = 'E7'; Value = '  +0.17%  ' },
    @{ Cell = 'E8'; Value = '  -0.26%  ' },
    @{ Cell = 'D9'; Value = '2.418.87' },
    @{ Cell = 'E9'; Value = '  -0.98%  ' },
    @{ Cell = 'E10'; Value = '  +1.34%  ' },
    @{ Cell = 'E11'; Value = '  +0.46%  ' },
    @{ Cell = 'D12'; Value = '5.12' },
    @{ Cell = 'E12'; Value = '  -1.80%  ' },
    @{ Cell = 'D13'; Value = '0.339' },
    @{ Cell = 'E13'; Value = '  -1.70%  ' },
    @{ Cell = 'D14'; Value = '26.11' },
    @{ Cell = 'E14'; Value = '  -1.30%  ' },
    @{ Cell = 'D15'; Value = '0.0000170' },
    @{ Cell = 'E15'; Value = '  -1.17%  ' },
    @{ Cell = 'D16'; Value = '2.828.38' },
    @{ Cell = 'E16'; Value = '  -1.53%  ' },
    @{ Cell = 'D17'; Value = '60.991.89' },
    @{ Cell = 'E17'; Value = '  -1.79%  ' },
    @{ Cell = 'D18'; Value = '2.424.78' },
    @{ Cell = 'E18'; Value = '  -0.73%  ' },
    @{ Cell = 'D19'; Value = '10.56' },
    @{ Cell = 'E19'; Value = '  -3.14%  ' },
    @{ Cell = 'D20'; Value = '7.28' },
    @{ Cell = 'E20'; Value = '  +2.14%  ' },
    @{ Cell = 'D21'; Value = '323.54' },
    @{ Cell = 'E21'; Value = '  -1.84%  ' },
    @{ Cell = 'D22'; Value = '4.04' },
    @{ Cell = 'E22'; Value = '  -1.55%  ' },
    @{ Cell = 'D23'; Value = '6.12' },
    @{ Cell = 'E23'; Value = '  +2.21%  ' },
    @{ Cell = 'E24'; Value = '  +0.05%  ' },
    @{ Cell = 'D25'; Value = '1.88' },
    @{ Cell = 'E25'; Value = '  -4.98%  ' },
    @{ Cell = 'D26'; Value = '64.83' },
    @{ Cell = 'D27'; Value = '8.81' },
    @{ Cell = 'E27'; Value = '  -6.18%  ' },
    @{ Cell = 'D28'; Value = '574.95' },
    @{ Cell = 'E28'; Value = '  -7.11%  ' },
    @{ Cell = 'D29'; Value = '2.569.17' },
    @{ Cell = 'E29'; Value = '  +0.05%  ' },
    @{ Cell = 'D30'; Value = '0.0₃0912' },
    @{ Cell = 'E30'; Value = '  -4.61%  ' },
    @{ Cell = 'D31'; Value = '7.87' },
    @{ Cell = 'E31'; Value = '  -1.73%  ' },
    @{ Cell = 'D32'; Value = '1.34' },
    @{ Cell = 'E32'; Value = '  -6.17%  ' },
    @{ Cell = 'D33'; Value = '1.83' },
    @{ Cell = 'E33'; Value = '  -2.66%  ' },
    @{ Cell = 'E34'; Value = '  -6.23%  ' },
    @{ Cell = 'E35'; Value = '  +0.24%  ' },
    @{ Cell = 'D36'; Value = '4.60' },
    @{ Cell = 'E36'; Value = '  -6.83%  ' },
    @{ Cell = 'D37'; Value = '0.366' },
    @{ Cell = 'E37'; Value = '  -3.14%  ' },
    @{ Cell = 'D38'; Value = '149.72' },
    @{ Cell = 'E38'; Value = '  -1.18%  ' },
    @{ Cell = 'E39'; Value = '  -3.78%  ' },
    @{ Cell = 'D40'; Value = '18.26' },
    @{ Cell = 'E40'; Value = '  -0.37%  ' },
    @{ Cell = 'D41'; Value = '5.10' },
    @{ Cell = 'E41'; Value = '  -3.00%  ' },
    @{ Cell = 'D43'; Value = '41.70' },
    @{ Cell = 'E43'; Value = '  -1.70%  ' },
    @{ Cell = 'D44'; Value = '1.65' },
    @{ Cell = 'E44'; Value = '  -6.24%  ' },
    @{ Cell = 'D45'; Value = '2.34' },
    @{ Cell = 'E45'; Value = '  -5.58%  ' },
    @{ Cell = 'E46'; Value = '  +20.19%  ' },
    @{ Cell = 'D47'; Value = '141.06' },
    @{ Cell = 'E47'; Value = '  -1.65%  ' },
    @{ Cell = 'D48'; Value = '3.51' },
    @{ Cell = 'E48'; Value = '  -3.43%  ' },
    @{ Cell = 'D49'; Value = '0.594' },
    @{ Cell = 'E49'; Value = '  -0.71%  ' },
    @{ Cell = 'D50'; Value = '19.49' },
    @{ Cell = 'E50'; Value = '  +0.07%  ' },
    @{ Cell = 'D51'; Value = '0.0505' },
    @{ Cell = 'E51'; Value = '  -3.84%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = "'" + $u.Value
    $range.Style = "Normal"
}
